# Applies the LOQ4053.docx content-reorganization edit described in the commit diff.
# Strategy: the edit is a pure reshuffle of text blocks among a fixed set of paragraphs
# (the paragraph count and styles stay the same); we overwrite each affected paragraph
# Range with its target OOXML using Range.InsertXML, which replaces just that range's
# contents in place (paragraph mark / neighboring paragraphs are left untouched).

$d = $word.ActiveDocument

$p7Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>Introduction to data calculations in Chemical Engineering; materials balance; Energy balances; materials and energy balances combined; Balances in processes in transient state.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(7).Range.InsertXML($p7Xml)

$p9Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Introdução aos cálculos em Engenharia Química; Balanços materiais; Balanços de energia; Balanços material e energético combinados; Balanços em processos no estado transiente.</w:t><w:br/></w:r><w:r><w:t>Introduzir e estabelecer aos alunos os princípios básicos a serem utilizados em todos os processos que envolvam a conservação de massa e energia. Esta disciplina propicia a realização de balanços globais de massa e energia em diferentes processos químicos evidenciando a importância da aplicação desta metodologia no projeto e otimização de processos químicos industriais.</w:t><w:br/></w:r><w:r><w:t>1 - Balanços Materiais</w:t><w:br/><w:t>- Introdução aos Balanços Materiais</w:t><w:br/><w:t>- Balanços Materiais que não envolvem reações químicas.</w:t><w:br/><w:t>- Balanços Materiais envolvendo reações químicas</w:t><w:br/><w:t>- Balanços Materiais com recirculação (reciclo e Bypass).</w:t><w:br/><w:t>2 - Balanços de Energia</w:t><w:br/><w:t>- Definições e conceitos. Formas de energia, calor, entalpia, valores de entalpia e capacidade calorífica.</w:t><w:br/><w:t>- Balanços de energia que não envolvem reações químicas.</w:t><w:br/><w:t>- Balanços de energia envolvendo reações químicas.</w:t><w:br/><w:t>3 - Balanços de Massa e de Energia Combinados</w:t><w:br/><w:t>- Aplicação combinada dos balanços de massa e energia em processos tais como umidificação, dissolução, processos de mistura, etc.</w:t><w:br/><w:t>4 - Balanços de Massa e de Energia em Processos de Regime Transiente</w:t><w:br/></w:r><w:r><w:t>Provas escritas; -participação e conteúdo de trabalho e seminário;</w:t><w:br/></w:r><w:r><w:t>Média Final = (Prova1 + 2xProva2 + Nota de Trabalho) / 4</w:t><w:br/><w:t>Média final mínima de aprovação = 5,0</w:t><w:br/></w:r><w:r><w:t>(Prova escrita + Média Final)/2</w:t><w:br/><w:t>Nota Final mínima para aprovação= 5,0</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(9).Range.InsertXML($p9Xml)

$p11Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>HIMMELBLAU, David M.  Eng. Química princípios e cálculos.  7. ed. LTC Editora,2006.</w:t><w:br/><w:t xml:space="preserve">GOMIDE, R.  Estequiometria Industrial. 3.ed. São Paulo: Ed. do Autor, 1984. </w:t><w:br/><w:t>FELDER, R.M; ROUSSEAU, R.W. Princípios elementares dos processos químicos. 3. ed. Rio de Janeiro: LTC Editora, 2005</w:t><w:br/><w:t xml:space="preserve">MOUYEN, O.A.; WATSON, K. M. AND RAGATZ, R.A.  Princípios dos processos químicos.  Lisboa: Lopes da Silva Editora, 2005. v. 1 </w:t><w:br/><w:t>CREMASCO, M. A. Fundamentos de transferência de massa. 1.ed. Campinas: Editora da UNICAMP, 1998.</w:t><w:br/><w:t>BRASIL, N. I. Introdução à Engenharia Química.  2. ed. Rio de Janeiro: Editora Interciência , 2004.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(11).Range.InsertXML($p11Xml)

$p12Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>And introduce students to the basic setting to be used in all cases involving the conservation of mass and energy principles. This course provides the realization of global balances of mass and energy in different chemical processes highlighting the importance of using this methodology in the design and optimization of industrial chemical processes.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(12).Range.InsertXML($p12Xml)

$p14Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>O grupo social alvo da atividade será composto pelos funcionários de indústrias que atuam em processos físicos, físico-químicos, biotecnológicos, envolvendo ou não questões ambientais, nos quais haja possibilidade de inserção de otimização do ponto de vista de balanço material e de energia. Este grupo inclui profissionais que trabalham diretamente na operação e gestão desses processos industriais, como engenheiros de processos, técnicos de operação, gestores de produção, entre outros.</w:t><w:br/><w:t>Esses profissionais são fundamentais para o funcionamento e a eficiência das indústrias, pois estão envolvidos na implementação, monitoramento e melhoria contínua dos processos produtivos. Eles têm um papel-chave na identificação de oportunidades de otimização, redução de custos, aumento da eficiência e minimização dos impactos ambientais associados às operações industriais.</w:t><w:br/><w:t>Portanto, ao realizar atividades de extensão com esse grupo social, o objetivo é capacitá-los e engajá-los em processos de aprendizado colaborativo, troca de conhecimentos e desenvolvimento de soluções práticas para os desafios enfrentados pela indústria. Essas atividades podem incluir análises de processos, desenvolvimento 3 de soluções de otimização, implementação de melhorias e avaliação de resultados, visando contribuir para o desenvolvimento profissional dos funcionários e para o aumento da competitividade e sustentabilidade das indústrias.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(14).Range.InsertXML($p14Xml)

$p17Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Método: </w:t></w:r><w:r><w:t>A realização de uma atividade extensionista voltada para funcionários de indústrias que atuam em processos físicos, físico-químicos, biotecnológicos envolvendo ou não questões ambientais, com possibilidade de inserção de otimização do ponto de vista de balanço material e de energia, pode trazer diversos benefícios tanto para os discentes envolvidos quanto para o grupo social alvo da ação. A atividade extensionista provavelmente envolverá conhecimentos de diversas disciplinas importantes dos cursos de engenharia química, bioquímica e ambiental, permitindo aos alunos uma visão mais ampla e integrada dos desafios industriais. Participar de atividades de extensão pode motivar os funcionários, aumentar seu conhecimento e habilidades, e promover um maior engajamento com os objetivos e valores da empresa. Através da aplicação de métodos de otimização de balanço material e energético, os funcionários poderão identificar oportunidades de melhorias nos processos industriais, aumentando a eficiência e reduzindo custos. Ao otimizar os processos industriais, é possível reduzir o consumo de recursos naturais, minimizar a geração de resíduos e emissões poluentes, contribuindo para a sustentabilidade ambiental e atendimento às legislações ambientais. Processos mais eficientes e sustentáveis podem tornar as empresas mais competitivas no mercado, melhorando sua posição</w:t><w:br/><w:t>frente à concorrência e aumentando sua rentabilidade a longo prazo.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Critério: </w:t></w:r><w:r><w:t>alunos e os funcionários podem trabalhar juntos para realizar uma análise detalhada dos processos</w:t><w:br/><w:t>industriais em vigor na empresa. Isso incluiria a identificação de pontos de entrada e saída de matéria e energia, o mapeamento de fluxos de material e energia, e a avaliação da eficiência dos processos.</w:t><w:br/><w:t>Atividades dos alunos:</w:t><w:br/><w:t xml:space="preserve">• Realizar levantamentos de campo para coletar dados sobre os processos industriais. </w:t><w:br/><w:t>• Utilizar ferramentas de análise, como balanços materiais e energéticos, para avaliar a eficiência dos processos.</w:t><w:br/><w:t>• Desenvolver relatórios e apresentações para documentar os resultados da análise.</w:t><w:br/><w:t xml:space="preserve">Atividades dos Funcionários </w:t><w:br/><w:t>• Fornecer informações e insights sobre os processos industriais, incluindo desafios e oportunidades de melhoria.</w:t><w:br/><w:t>• Colaborar com os alunos na interpretação dos dados coletados e na identificação de áreas para otimização.</w:t><w:br/><w:t>• Participar de discussões e sessões de brainstorming para gerar ideias de melhoria dos processos.</w:t><w:br/><w:t>2. DESENVOLVIMENTO DE SOLUÇÕES DE OTIMIZAÇÃO</w:t><w:br/><w:t>Com base na análise dos processos industriais, alunos e funcionários podem colaborar no desenvolvimento de soluções de otimização para melhorar a eficiência dos processos, reduzir custos e minimizar impactos ambientais.</w:t><w:br/><w:t>Atividades dos alunos:</w:t><w:br/><w:t>• Pesquisar e identificar tecnologias e práticas de ponta para otimização de processos industriais.</w:t><w:br/><w:t>• Propor soluções inovadoras e viáveis para melhorar a eficiência e sustentabilidade dos processos.</w:t><w:br/><w:t>4</w:t><w:br/><w:t>• Desenvolver modelos de simulação para avaliar o desempenho das soluções propostas.</w:t><w:br/><w:t>Atividades dos Funcionários</w:t><w:br/><w:t>• Contribuir com conhecimento prático e experiência na avaliação da viabilidade e aplicabilidade das soluções propostas.</w:t><w:br/><w:t>• Participar de sessões de revisão e refinamento das soluções propostas.</w:t><w:br/><w:t>• Identificar possíveis obstáculos ou desafios na implementação das soluções e colaborar na busca de estratégias para superá-los.</w:t><w:br/><w:t>Essas atividades propostas visam promover uma colaboração eficaz entre alunos e funcionários, proporcionando uma oportunidade única de aprendizado prático e aplicado, ao mesmo tempo em que busca soluções reais para desafios enfrentados pela indústria.</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Norma de recuperação: </w:t></w:r><w:r><w:t>Para garantir que o grupo social externo à universidade possa avaliar adequadamente a atividade realizada em conjunto com os estudantes, durante sua realização e ao final, é importante estabelecer mecanismos de avaliação claros e transparentes.</w:t><w:br/><w:t>DURANTE A REALIZAÇÃO DA ATIVIDADE</w:t><w:br/><w:t>Comunicação Aberta e Transparente:</w:t><w:br/><w:t>• Manter canais de comunicação abertos entre os membros da equipe da indústria e os estudantes.</w:t><w:br/><w:t>• Realizar reuniões regulares para discutir o progresso, desafios e oportunidades de melhoria.</w:t><w:br/><w:t>• Permitir que os funcionários da indústria forneçam feedback contínuo sobre o desempenho e a</w:t><w:br/><w:t>contribuição dos alunos.</w:t><w:br/><w:t>Acompanhamento de Indicadores Chave:</w:t><w:br/><w:t>• Estabelecer indicadores de desempenho relevantes para a atividade, como eficiência dos processos, redução de custos ou impactos ambientais.</w:t><w:br/><w:t>• Monitorar e registrar regularmente os dados relacionados aos indicadores-chave para avaliar o</w:t><w:br/><w:t>progresso e identificar áreas para ajustes.</w:t><w:br/><w:t>Participação Ativa dos Funcionários:</w:t><w:br/><w:t>• Envolver os funcionários da indústria na supervisão e avaliação das atividades dos estudantes.</w:t><w:br/><w:t>• Solicitar que os funcionários forneçam comentários e avaliações específicas sobre o trabalho dos alunos, destacando pontos fortes e áreas de melhoria.</w:t><w:br/><w:t>AO FINAL DA ATIVIDADE</w:t><w:br/><w:t>Avaliação Conjunta dos Resultados:</w:t><w:br/><w:t>• Realizar uma avaliação conjunta dos resultados alcançados durante a atividade, incluindo benefícios obtidos, lições aprendidas e desafios enfrentados.</w:t><w:br/><w:t>• Conduzir sessões de discussão para revisar os resultados e identificar oportunidades de aplicação prática.</w:t><w:br/><w:t>Coleta de Feedback:</w:t><w:br/><w:t>• Solicitar feedback dos funcionários da indústria sobre a experiência de trabalhar com os estudantes.</w:t><w:br/><w:t>• Realizar entrevistas ou pesquisas de satisfação para coletar feedback sobre a qualidade do trabalho dos alunos, a eficácia das soluções propostas e o impacto percebido nas operações da indústria.</w:t><w:br/><w:t>Avaliação dos Objetivos Alcançados:</w:t><w:br/><w:t>• Comparar os resultados alcançados com os objetivos inicialmente estabelecidos para a atividade.</w:t><w:br/><w:t>5</w:t><w:br/><w:t>• Identificar lacunas entre as expectativas e os resultados reais e discutir estratégias para abordar essas lacunas no futuro.</w:t><w:br/><w:t>Ao seguir a avaliação apresentada, o grupo social externo à universidade poderá avaliar de forma abrangente e significativa a atividade realizada em conjunto com os estudantes, garantindo uma compreensão clara dos benefícios e impactos da colaboração</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(17).Range.InsertXML($p17Xml)

$p19Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>5817045 - Elisângela de Jesus Cândido Moraes</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(19).Range.InsertXML($p19Xml)

